$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 511; all rows 511:631 shift down to 512:632.
$ws.Rows.Item(511).Insert()

# Populate the new row 511 with the new weekly price-report record.
$ws.Range("A511").Value = 4
$ws.Range("B511").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C511").Value = "Los Lagos"
$ws.Range("D511").Value = 45173
$ws.Range("E511").Value = 10
$ws.Range("F511").Value = 100114013
$ws.Range("G511").Value = "Zanahoria"
$ws.Range("H511").Value = "Sin especificar"
$ws.Range("I511").Value = "Primera"
$ws.Range("J511").Value = 150
$ws.Range("K511").Value = 7000
$ws.Range("L511").Value = 7000
$ws.Range("M511").Value = 7000
$ws.Range("N511").Value = "$/saco 20 kilos"
$ws.Range("O511").Value = "Provincia de Llanquihue"
$ws.Range("P511").Value = 350
$ws.Range("Q511").Value = 20
$ws.Range("R511").Value = "Hortaliza"
